# Fix a typo in the Dutch hairnet instructions text on the "eeg" sheet.
# Cell C26 contained: "... en er modieus uit blijft zijn in het openbaar."
# It should read:      "... en er modieus uit blijft zien in het openbaar."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("eeg")

$enDash = [char]0x2013
$ws.Range("C26").Value = " zodat u uw haar kan bedekken met een beschermend haarnetje " + $enDash + " en er modieus uit blijft zien in het openbaar."
